$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Subjects section (column C): Desc -> FullName, MaxAbsencesAllowed -> MinAttendanceRequired
$ws.Range("C6").Value = "FullName"
$ws.Range("C9").Value = "MinAttendanceRequired"

# --- Enrollments section (column K): insert SubjectFK field, pushing the rest down one row
$ws.Range("K6").Value = "SubjectFK"
$ws.Range("K7").Value = "SubjectCode"
$ws.Range("K8").Value = "SubjectName"
$ws.Range("K9").Value = "Grade"
$ws.Range("K10").Value = "Shift"
$ws.Range("K11").Value = "StudentFK"

# --- Column C widened to fit the new "FullName" header text
$ws.Columns.Item(3).ColumnWidth = 21.666666666666664

# --- Restore the active selection to F4
$ws.Range("F4").Select()
